# #5: property boat&car done
# Rebuild the "car" sheet (汽車, 3rd sheet) header row with the full set of
# standard property-record field names (matching the 土地/建物 sheets) and
# extend the data row with the matching legislator/source metadata columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# --- Header row (row 1): full field-name header, same pattern used on the
#     land/building sheets, minus the land-only share_portion/portion/total
#     columns and with "capacity" standing in for "area".
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "capacity"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "register_date"
$ws.Range("F1").Value = "register_reason"
$ws.Range("G1").Value = "acquire_value"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# --- Data row (row 2): keep the existing car record values, and append the
#     property_category/category/date/legislator/source/index columns.
$ws.Range("H2").Value = "land"
$ws.Range("I2").Value = "normal"
$ws.Range("J2").Value = "2013-12-26"
$ws.Range("K2").Value = "林郁方"
$ws.Range("L2").Value = 716
$ws.Range("M2").Value = "tmp4c8a1"
$ws.Range("N2").Value = 50
